$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching style of existing header G1
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats: copy only formatting, keep the value
$excel.CutCopyMode = $false

# Add values for the new Save column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
